# chore: update Sheets via scheduled runner
# Refresh cached Universalis market-board price/profit figures for a batch
# of Leve rows across the per-job workbook tabs (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H:N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3127.3513
$ws.Range("I138").Value = 3267.3333
$ws.Range("J138").Value = 2994.7368
$ws.Range("K138").Value = 9801.999899999999
$ws.Range("L138").Value = 8984.2104
$ws.Range("M138").Value = -4661.999899999999
$ws.Range("N138").Value = -19264.2104

$ws.Range("H141").Value = 2548.4167
$ws.Range("I141").Value = 1262.4706
$ws.Range("J141").Value = 5671.4287
$ws.Range("K141").Value = 3787.4118
$ws.Range("L141").Value = 17014.2861
$ws.Range("M141").Value = 1392.5882
$ws.Range("N141").Value = -27374.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8777496
$ws.Range("I32").Value = 5604.727
$ws.Range("J32").Value = 250004500
$ws.Range("K32").Value = 5604.727
$ws.Range("L32").Value = 250004500
$ws.Range("M32").Value = -5317.727
$ws.Range("N32").Value = -250005074

$ws.Range("H61").Value = 3945823.2
$ws.Range("I61").Value = 3631253
$ws.Range("J61").Value = 4526568.5
$ws.Range("K61").Value = 3631253
$ws.Range("L61").Value = 4526568.5
$ws.Range("M61").Value = -3631041
$ws.Range("N61").Value = -4526992.5

$ws.Range("H97").Value = 393.64285
$ws.Range("I97").Value = 359.25
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 359.25
$ws.Range("L97").Value = 600
$ws.Range("M97").Value = 136.75
$ws.Range("N97").Value = -1592

$ws.Range("H122").Value = 2268.3794
$ws.Range("I122").Value = 2216.9524
$ws.Range("J122").Value = 2403.375
$ws.Range("K122").Value = 6650.8572
$ws.Range("L122").Value = 7210.125
$ws.Range("M122").Value = -4200.8572
$ws.Range("N122").Value = -12110.125

$ws.Range("H136").Value = 3945823.2
$ws.Range("I136").Value = 3631253
$ws.Range("J136").Value = 4526568.5
$ws.Range("K136").Value = 10893759
$ws.Range("L136").Value = 13579705.5
$ws.Range("M136").Value = -10891209
$ws.Range("N136").Value = -13584805.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 107253560
$ws.Range("I20").Value = 62525788
$ws.Range("J20").Value = 166890600
$ws.Range("K20").Value = 62525788
$ws.Range("L20").Value = 166890600
$ws.Range("M20").Value = -62525541
$ws.Range("N20").Value = -166891094

$ws.Range("H94").Value = 1735.7646
$ws.Range("I94").Value = 1043.4286
$ws.Range("J94").Value = 4966.6665
$ws.Range("K94").Value = 1043.4286
$ws.Range("L94").Value = 4966.6665
$ws.Range("M94").Value = -592.4286
$ws.Range("N94").Value = -5868.6665

$ws.Range("H99").Value = 1824.5
$ws.Range("I99").Value = 1795.25
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1795.25
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -297.25
$ws.Range("N99").Value = -4996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13159389
$ws.Range("I31").Value = 22728306
$ws.Range("J31").Value = 2127.4688
$ws.Range("K31").Value = 22728306
$ws.Range("L31").Value = 2127.4688
$ws.Range("M31").Value = -22728011
$ws.Range("N31").Value = -2717.4688

$ws.Range("H34").Value = 13159389
$ws.Range("I34").Value = 22728306
$ws.Range("J34").Value = 2127.4688
$ws.Range("K34").Value = 22728306
$ws.Range("L34").Value = 2127.4688
$ws.Range("M34").Value = -22728104
$ws.Range("N34").Value = -2531.4688

$ws.Range("H59").Value = 22975
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 22975
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 22975
$ws.Range("N59").Value = -25265

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 23.846153
$ws.Range("I33").Value = 30
$ws.Range("J33").Value = 20
$ws.Range("K33").Value = 180
$ws.Range("L33").Value = 120
$ws.Range("M33").Value = 103
$ws.Range("N33").Value = -686

$ws.Range("H34").Value = 1256.5
$ws.Range("I34").Value = 661
$ws.Range("J34").Value = 2447.5
$ws.Range("K34").Value = 1983
$ws.Range("L34").Value = 7342.5
$ws.Range("M34").Value = -1899
$ws.Range("N34").Value = -7510.5

$ws.Range("H39").Value = 2611.25
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2611.25
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7833.75
$ws.Range("N39").Value = -8421.75

$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = $null

$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -32122

$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -100608

$ws.Range("H93").Value = 6216.7334
$ws.Range("I93").Value = 1924
$ws.Range("J93").Value = 6523.357
$ws.Range("K93").Value = 5772
$ws.Range("L93").Value = 19570.071
$ws.Range("M93").Value = -3900
$ws.Range("N93").Value = -23314.071

$ws.Range("H131").Value = 7953394
$ws.Range("I131").Value = 250000160
$ws.Range("J131").Value = 17434.426
$ws.Range("K131").Value = 750000480
$ws.Range("L131").Value = 52303.278
$ws.Range("M131").Value = -749995440
$ws.Range("N131").Value = -62383.278

$ws.Range("H137").Value = 3669.4666
$ws.Range("I137").Value = 2888.75
$ws.Range("J137").Value = 4561.7144
$ws.Range("K137").Value = 8666.25
$ws.Range("L137").Value = 13685.1432
$ws.Range("M137").Value = -3566.25
$ws.Range("N137").Value = -23885.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 13159297
$ws.Range("I97").Value = 919.2308
$ws.Range("J97").Value = 41669116
$ws.Range("K97").Value = 919.2308
$ws.Range("L97").Value = 41669116
$ws.Range("M97").Value = -423.2308
$ws.Range("N97").Value = -41670108

$ws.Range("H107").Value = 235.5
$ws.Range("I107").Value = 164
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 164
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1756
$ws.Range("N107").Value = -4290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 84581.75
$ws.Range("I16").Value = 84581.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 84581.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -84411.75
$ws.Range("N16").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 29500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 29500
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -30130

$ws.Range("H73").Value = 29500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 29500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 29500
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -31684

$ws.Range("H87").Value = 32000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 32000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 32000
$ws.Range("N87").Value = -34496

$ws.Range("H90").Value = 32000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 32000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 96000
$ws.Range("N90").Value = -108480

$ws.Range("H107").Value = 28842.715
$ws.Range("I107").Value = 40199.8
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 120599.4
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = -118679.4
$ws.Range("N107").Value = -5190
